$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 21, shifting existing rows 21-48 down to 23-50.
$ws.Rows.Item(21).Resize(2).Insert()

# Row 21: new "Especial" entry for date 2022-03-02
$ws.Cells.Item(21, 1).Value = 2
$ws.Cells.Item(21, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(21, 3).Value = "Coquimbo"
$ws.Cells.Item(21, 4).Value = 44622
$ws.Cells.Item(21, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100107
$ws.Cells.Item(21, 8).Value = "Otros"
$ws.Cells.Item(21, 9).Value = 100107011
$ws.Cells.Item(21, 10).Value = "Tuna"
$ws.Cells.Item(21, 11).Value = "Sin especificar"
$ws.Cells.Item(21, 12).Value = "Especial"
$ws.Cells.Item(21, 13).Value = 240
$ws.Cells.Item(21, 14).Value = 13000
$ws.Cells.Item(21, 15).Value = 14000
$ws.Cells.Item(21, 16).Value = 13500
$ws.Cells.Item(21, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(21, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(21, 19).Value = 750
$ws.Cells.Item(21, 20).Value = 18

# Row 22: new "Primera" entry for date 2022-03-02
$ws.Cells.Item(22, 1).Value = 2
$ws.Cells.Item(22, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value = 44622
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(22, 5).Value = 4
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100107
$ws.Cells.Item(22, 8).Value = "Otros"
$ws.Cells.Item(22, 9).Value = 100107011
$ws.Cells.Item(22, 10).Value = "Tuna"
$ws.Cells.Item(22, 11).Value = "Sin especificar"
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 500
$ws.Cells.Item(22, 14).Value = 11000
$ws.Cells.Item(22, 15).Value = 12000
$ws.Cells.Item(22, 16).Value = 11500
$ws.Cells.Item(22, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(22, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 19).Value = 639
$ws.Cells.Item(22, 20).Value = 18
